# Add another formula for operators
$wb = $excel.ActiveWorkbook

# Work on the "Final" worksheet (sheet2.xml) — add a new formula row.
$ws = $wb.Worksheets.Item("Final")

# New row 13, column B: formula for operator precedence practice.
$ws.Range("B13").Formula = "=13-2+(15/3)^2"

# Make "Final" the active/selected sheet with the active cell at G12,
# matching the author's recorded selection state.
$ws.Activate() | Out-Null
$ws.Range("G12").Select() | Out-Null

$wb.Save() | Out-Null
